# StateList.xlsx edit: refresh the state-adjacency ("M-Colouring") data
# and add the new "Ladakh" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the adjacency list (column C) for states whose neighbour
#     list changed (values re-sorted / corrected, and the new state
#     "Ladakh" (#37) added as a neighbour of J&K and Himachal Pradesh). ---
$ws.Range("C2").Value  = "2,3,37"
$ws.Range("C3").Value  = "1,3,5,6,37"
$ws.Range("C6").Value  = "2,6,9"
$ws.Range("C7").Value  = "2,3,4,5,7,8,9"
$ws.Range("C10").Value = "5,6,7,8,10,20,22,23"
$ws.Range("C21").Value = "9,10,19,21,22,23"
$ws.Range("C23").Value = "9,20,21,23,27,28,36"
$ws.Range("C24").Value = "8,9,22,24,27"
$ws.Range("C25").Value = "8,23,25,26,27"
$ws.Range("C28").Value = "22,23,24,26,29,30,36"
$ws.Range("C29").Value = "21,22,29,33,34,36"
$ws.Range("C30").Value = "27,28,30,32,33,36"
$ws.Range("C34").Value = "28,29,32,34"
$ws.Range("C35").Value = "28,32,33"
$ws.Range("C37").Value = "22,27,28,29"

# --- Add the new state "Ladakh" as row 38 (state code 37). ---
$ws.Range("A38").Value = 37
$ws.Range("B38").Value = "Ladakh"
$ws.Range("C38").Value = "1,2"

# --- Update the view: scroll down, zoom in a bit, and move the
#     selection to C30 (matches the author's saved UI state). ---
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
$win.Zoom = 119
$ws.Range("C30").Select()
